$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.873.68'
$ws.Range("E2").Value = '  -5.91%  '

$ws.Range("D3").Value = '3.100.84'
$ws.Range("E3").Value = '  -6.28%  '

$ws.Range("E4").Value = '  +0.23%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '515.34'
$ws.Range("E5").Value = '  -7.77%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '130.93'
$ws.Range("E6").Value = '  -8.71%  '

$ws.Range("E7").Value = '  +0.10%  '

$ws.Range("D8").Value = '3.098.02'
$ws.Range("E8").Value = '  -6.48%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.441'
$ws.Range("E9").Value = '  -7.54%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.05'
$ws.Range("E10").Value = '  -10.24%  '

$ws.Range("E11").Value = '  -10.79%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.375'
$ws.Range("E12").Value = '  -8.35%  '

$ws.Range("D13").Value = '3.671.75'
$ws.Range("E13").Value = '  -5.10%  '

$ws.Range("E14").Value = '  -2.44%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '25.12'
$ws.Range("E15").Value = '  -7.77%  '

$ws.Range("D16").Value = '3.117.53'
$ws.Range("E16").Value = '  -5.30%  '

$ws.Range("D17").Value = '57.048.25'
$ws.Range("E17").Value = '  -5.53%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0000147'
$ws.Range("E18").Value = '  -11.76%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.67'
$ws.Range("E19").Value = '  -7.84%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.72'
$ws.Range("E20").Value = '  -11.50%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.82'
$ws.Range("E21").Value = '  -9.03%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '340.92'
$ws.Range("E22").Value = '  -9.36%  '

$ws.Range("E23").Value = '  -0.01%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '67.79'
$ws.Range("E24").Value = '  -8.34%  '

$ws.Range("E25").Value = '  -9.09%  '

$ws.Range("D26").Value = '3.258.44'
$ws.Range("E26").Value = '  -5.19%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.997'
$ws.Range("E27").Value = '  +0.25%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.160'
$ws.Range("E28").Value = '  -7.26%  '

$ws.Range("D29").Value = '0.0₃0916'
$ws.Range("E29").Value = '  -12.04%  '

$ws.Range("E30").Value = '  -0.06%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.57'
$ws.Range("E31").Value = '  -9.61%  '

$ws.Range("E32").Value = '  -11.61%  '

$ws.Range("E33").Value = '  -10.74%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '21.31'
$ws.Range("E34").Value = '  -5.60%  '

$ws.Range("E35").Value = '  -7.43%  '

$ws.Range("B36").Value = 'Monero'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '155.93'
$ws.Range("E36").Value = '  -6.37%  '

$ws.Range("B37").Value = 'NEARProtocol'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.72'
$ws.Range("E37").Value = '  -9.91%  '

$ws.Range("E38").Value = '  -9.94%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.36'
$ws.Range("E39").Value = '  -11.66%  '

$ws.Range("D41").Value = '3.166.73'
$ws.Range("E41").Value = '  -4.73%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0676'
$ws.Range("E42").Value = '  -9.20%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '40.09'
$ws.Range("E43").Value = '  -4.40%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.676'
$ws.Range("E44").Value = '  -10.04%  '

$ws.Range("E45").Value = '  +0.09%  '

$ws.Range("E46").Value = '  -8.98%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.03'
$ws.Range("E47").Value = '  -8.60%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.44'
$ws.Range("E48").Value = '  -10.26%  '

$ws.Range("D49").Value = '2.222.23'
$ws.Range("E49").Value = '  -5.75%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.06'
$ws.Range("E50").Value = '  -7.31%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '19.64'
$ws.Range("E51").Value = '  -7.63%  '

